$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4.918333463945449
$ws.Range("D2").Value = 9.871317234904735
$ws.Range("E2").Value = 13.94456563685932
$ws.Range("F2").Value = 29.37012180912973
$ws.Range("G2").Value = 28.74752955850481
$ws.Range("H2").Value = 14.20720873030578
$ws.Range("J2").Value = 9.778814238681491
$ws.Range("K2").Value = 13.88362409779512
$ws.Range("N2").Value = 16.77597527619096
$ws.Range("O2").Value = 21.64094833031352
$ws.Range("C3").Value = 4.74976546125179
$ws.Range("D3").Value = 9.819010985222608
$ws.Range("E3").Value = 13.8829927839413
$ws.Range("F3").Value = 29.38081368111657
$ws.Range("G3").Value = 28.73608409569469
$ws.Range("H3").Value = 14.25218378735418
$ws.Range("J3").Value = 9.784284122349508
$ws.Range("K3").Value = 13.38930414262038
$ws.Range("N3").Value = 16.79693341726779
$ws.Range("O3").Value = 21.69849498478166
$ws.Range("C4").Value = 4.644625336893507
$ws.Range("D4").Value = 9.788375328612611
$ws.Range("E4").Value = 13.84794121546453
$ws.Range("F4").Value = 29.3956068296843
$ws.Range("G4").Value = 28.73998121531753
$ws.Range("H4").Value = 14.28249324912956
$ws.Range("J4").Value = 9.789240339436391
$ws.Range("K4").Value = 13.07717360696458
$ws.Range("N4").Value = 16.81154983180017
$ws.Range("O4").Value = 21.73939087741977
$ws.Range("C5").Value = 4.601445123842196
$ws.Range("D5").Value = 9.776272994533617
$ws.Range("E5").Value = 13.83436115084144
$ws.Range("F5").Value = 29.40370173265752
$ws.Range("G5").Value = 28.7443112998841
$ws.Range("H5").Value = 14.29552119383865
$ws.Range("J5").Value = 9.791662118937493
$ws.Range("K5").Value = 12.94799334580957
$ws.Range("N5").Value = 16.81794655352222
$ws.Range("O5").Value = 21.75745030022167
$ws.Range("C6").Value = 4.594257241839066
$ws.Range("D6").Value = 9.774286751267034
$ws.Range("E6").Value = 13.83214900219099
$ws.Range("F6").Value = 29.4051706025164
$ws.Range("G6").Value = 28.7451957106395
$ws.Range("H6").Value = 14.29772530935384
$ws.Range("J6").Value = 9.792088545565965
$ws.Range("K6").Value = 12.9264288777339
$ws.Range("N6").Value = 16.819035349034
$ws.Range("O6").Value = 21.76053310812478
$ws.Range("C7").Value = 4.644044242228647
$ws.Range("D7").Value = 9.788210553468478
$ws.Range("E7").Value = 13.84775520655808
$ws.Range("F7").Value = 29.39570763719441
$ws.Range("G7").Value = 28.7400285191
$ws.Range("H7").Value = 14.28266621038471
$ws.Range("J7").Value = 9.789271372044125
$ws.Range("K7").Value = 13.07543921820598
$ws.Range("N7").Value = 16.81163431588888
$ws.Range("O7").Value = 21.73962879513242
$ws.Range("C8").Value = 4.860595803317272
$ws.Range("D8").Value = 9.852980675890986
$ws.Range("E8").Value = 13.92277057757975
$ws.Range("F8").Value = 29.37209946871038
$ws.Range("G8").Value = 28.74131352105405
$ws.Range("H8").Value = 14.22215628357825
$ws.Range("J8").Value = 9.780368885320991
$ws.Range("K8").Value = 13.71507448195352
$ws.Range("N8").Value = 16.78283935677976
$ws.Range("O8").Value = 21.65963334715347
$ws.Range("C9").Value = 5.26916979962613
$ws.Range("D9").Value = 9.991297028731976
$ws.Range("E9").Value = 14.09120175203906
$ws.Range("F9").Value = 29.39114715638016
$ws.Range("G9").Value = 28.8306118476243
$ws.Range("H9").Value = 14.12492008446184
$ws.Range("J9").Value = 9.775569798112274
$ws.Range("K9").Value = 14.89394086934353
$ws.Range("N9").Value = 16.74020506714318
$ws.Range("O9").Value = 21.54708745055604
$ws.Range("C10").Value = 5.555865203161962
$ws.Range("D10").Value = 10.09918227045378
$ws.Range("E10").Value = 14.22718562818184
$ws.Range("F10").Value = 29.44496955406417
$ws.Range("G10").Value = 28.9490590387306
$ws.Range("H10").Value = 14.06659879043486
$ws.Range("J10").Value = 9.779730785154159
$ws.Range("K10").Value = 15.70576382677727
$ws.Range("N10").Value = 16.71726468405396
$ws.Range("O10").Value = 21.49166584878596
$ws.Range("C11").Value = 5.682709413931661
$ws.Range("D11").Value = 10.14948087416241
$ws.Range("E11").Value = 14.29152686505045
$ws.Range("F11").Value = 29.47807216704331
$ws.Range("G11").Value = 29.01434057786903
$ws.Range("H11").Value = 14.04292717110616
$ws.Range("J11").Value = 9.783283859823239
$ws.Range("K11").Value = 16.06186263797138
$ws.Range("N11").Value = 16.70863797121444
$ws.Range("O11").Value = 21.47241840695168
$ws.Range("C12").Value = 5.730181488997246
$ws.Range("D12").Value = 10.16869081508328
$ws.Range("E12").Value = 14.31623194872603
$ws.Range("F12").Value = 29.49184150698688
$ws.Range("G12").Value = 29.04068884668659
$ws.Range("H12").Value = 14.03437539689752
$ws.Range("J12").Value = 9.784867039928537
$ws.Range("K12").Value = 16.19470925681868
$ws.Range("N12").Value = 16.70563043816046
$ws.Range("O12").Value = 21.4659905316413
$ws.Range("C13").Value = 5.719983172509486
$ws.Range("D13").Value = 10.16454656713659
$ws.Range("E13").Value = 14.31089639215032
$ws.Range("F13").Value = 29.48882124550105
$ws.Range("G13").Value = 29.03494210668707
$ws.Range("H13").Value = 14.03619882788569
$ws.Range("J13").Value = 9.784515519509524
$ws.Range("K13").Value = 16.16618880375291
$ws.Range("N13").Value = 16.70626665030904
$ws.Range("O13").Value = 21.46733656471435
$ws.Range("C14").Value = 5.686626464646069
$ws.Range("D14").Value = 10.15105809504333
$ws.Range("E14").Value = 14.29355263985054
$ws.Range("F14").Value = 29.47918028100154
$ws.Range("G14").Value = 29.01647572439873
$ws.Range("H14").Value = 14.04221534556956
$ws.Range("J14").Value = 9.783409351372116
$ws.Range("K14").Value = 16.07283263366845
$ws.Range("N14").Value = 16.70838534990208
$ws.Range("O14").Value = 21.47187231570479
$ws.Range("C15").Value = 5.666120148225583
$ws.Range("D15").Value = 10.1428168581515
$ws.Range("E15").Value = 14.28297293173751
$ws.Range("F15").Value = 29.47343544729622
$ws.Range("G15").Value = 29.00537609197508
$ws.Range("H15").Value = 14.04595434083096
$ws.Range("J15").Value = 9.78276271511254
$ws.Range("K15").Value = 16.01538599556867
$ws.Range("N15").Value = 16.70971684463417
$ws.Range("O15").Value = 21.4747627684468
$ws.Range("C16").Value = 5.547499554277927
$ws.Range("D16").Value = 10.09591863511573
$ws.Range("E16").Value = 14.22302932421561
$ws.Range("F16").Value = 29.44297923365997
$ws.Range("G16").Value = 28.94502119441058
$ws.Range("H16").Value = 14.06820340639011
$ws.Range("J16").Value = 9.779531906359397
$ws.Range("K16").Value = 15.68221725578156
$ws.Range("N16").Value = 16.7178647837984
$ws.Range("O16").Value = 21.49304400083231
$ws.Range("C17").Value = 5.473778818236203
$ws.Range("D17").Value = 10.06745193233092
$ws.Range("E17").Value = 14.18687996810518
$ws.Range("F17").Value = 29.42649949187375
$ws.Range("G17").Value = 28.9109077440955
$ws.Range("H17").Value = 14.08258549957627
$ws.Range("J17").Value = 9.777974443109652
$ws.Range("K17").Value = 15.47437076180775
$ws.Range("N17").Value = 16.7233259615144
$ws.Range("O17").Value = 21.50578882664655
$ws.Range("C18").Value = 5.431042552639622
$ws.Range("D18").Value = 10.05119434476596
$ws.Range("E18").Value = 14.16632240479539
$ws.Range("F18").Value = 29.41783240504458
$ws.Range("G18").Value = 28.89236043936987
$ws.Range("H18").Value = 14.09112674565803
$ws.Range("J18").Value = 9.777234969616366
$ws.Range("K18").Value = 15.35358537246947
$ws.Range("N18").Value = 16.7266374125062
$ws.Range("O18").Value = 21.51368052730312
$ws.Range("C19").Value = 5.41651693738778
$ws.Range("D19").Value = 10.04571005845313
$ws.Range("E19").Value = 14.15940274747534
$ws.Range("F19").Value = 29.41503741099116
$ws.Range("G19").Value = 28.88626538020917
$ws.Range("H19").Value = 14.09406484176621
$ws.Range("J19").Value = 9.777011476957908
$ws.Range("K19").Value = 15.31248027208473
$ws.Range("N19").Value = 16.72778789404686
$ws.Range("O19").Value = 21.51644880976215
$ws.Range("C20").Value = 5.481661447123098
$ws.Range("D20").Value = 10.07047037692385
$ws.Range("E20").Value = 14.19070396635703
$ws.Range("F20").Value = 29.42816982160101
$ws.Range("G20").Value = 28.91442811434776
$ws.Range("H20").Value = 14.08102664829717
$ws.Range("J20").Value = 9.778124064341348
$ws.Range("K20").Value = 15.49662523881787
$ws.Range("N20").Value = 16.72272698826003
$ws.Range("O20").Value = 21.50437401095041
$ws.Range("C21").Value = 5.696439713162225
$ws.Range("D21").Value = 10.15501566575498
$ws.Range("E21").Value = 14.2986378139387
$ws.Range("F21").Value = 29.48197862195311
$ws.Range("G21").Value = 29.02185568462038
$ws.Range("H21").Value = 14.0404369547474
$ws.Range("J21").Value = 9.7837278173354
$ws.Range("K21").Value = 16.10030863618897
$ws.Range("N21").Value = 16.70775600882291
$ws.Range("O21").Value = 21.47051667378126
$ws.Range("C22").Value = 5.833521332322486
$ws.Range("D22").Value = 10.21121467788706
$ws.Range("E22").Value = 14.37115525463883
$ws.Range("F22").Value = 29.52433558413477
$ws.Range("G22").Value = 29.10154499021586
$ws.Range("H22").Value = 14.01631204036677
$ws.Range("J22").Value = 9.788775162179476
$ws.Range("K22").Value = 16.483153286549
$ws.Range("N22").Value = 16.69948207670062
$ws.Range("O22").Value = 21.4534069622016
$ws.Range("C23").Value = 5.760673311438646
$ws.Range("D23").Value = 10.18113808828272
$ws.Range("E23").Value = 14.3322760464695
$ws.Range("F23").Value = 29.5010731208806
$ws.Range("G23").Value = 29.05815058880165
$ws.Range("H23").Value = 14.02896777651786
$ws.Range("J23").Value = 9.785954950338571
$ws.Range("K23").Value = 16.2799224732079
$ws.Range("N23").Value = 16.70376012979949
$ws.Range("O23").Value = 21.46207868354761
$ws.Range("C24").Value = 5.478098807238517
$ws.Range("D24").Value = 10.06910539983439
$ws.Range("E24").Value = 14.18897443410212
$ws.Range("F24").Value = 29.42741215033235
$ws.Range("G24").Value = 28.91283323677049
$ws.Range("H24").Value = 14.08173055529735
$ws.Range("J24").Value = 9.77805593479197
$ws.Range("K24").Value = 15.48656800516409
$ws.Range("N24").Value = 16.72299724900945
$ws.Range("O24").Value = 21.50501189086146
$ws.Range("C25").Value = 5.160765861609627
$ws.Range("D25").Value = 9.952731771626739
$ws.Range("E25").Value = 14.04342918994351
$ws.Range("F25").Value = 29.37899389176707
$ws.Range("G25").Value = 28.79715856651402
$ws.Range("H25").Value = 14.14892541308796
$ws.Range("J25").Value = 9.775569798112274
$ws.Range("K25").Value = 14.58402414774255
$ws.Range("N25").Value = 16.75026295650373
$ws.Range("O25").Value = 21.57276214382952